$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary field updates ---
$ws.Range("D5").Value = "Report Generated On: 08/16/2025 12:47 AM"
$ws.Range("C8").Value = 16891.47
$ws.Range("C9").Value = 86
$ws.Range("C10").Value = "07/28/2025 to 08/03/25"
$ws.Range("G13").Value = "709-2"

# --- Append a new "Thursday (07/31/2025)" daily section after the existing
#     data (rows 1-107), mirroring the style pattern used by the other daily
#     sections (e.g. the "Wednesday" block at rows 80-86 and the closing
#     TOTAL row at row 107). ---

# Row 110 is a merged single-line day header (like A80:H80). Merge the
# (currently blank) destination row first so the merge doesn't propagate the
# anchor cell's style across the whole range, then paste in just the
# formatting from the analogous row, then set the text.
$ws.Range("A110:H110").Merge()
$ws.Range("A80").Copy()
$ws.Range("A110").PasteSpecial(-4122)
$ws.Range("A110").Value = "Thursday (07/31/2025)"

# Rows 111-116: column-header row + 5 data rows, copied wholesale from the
# analogous block (A81:H86) so the alternating row styles (s=8/9/10/11/12/13/14)
# come along for free.
$ws.Range("A81:H86").Copy($ws.Range("A111"))

# Row 117 is the merged TOTAL row (like A107:G107). Same merge-before-paste
# trick so only A117/H117 end up styled.
$ws.Range("A117:G117").Merge()
$ws.Range("A107").Copy()
$ws.Range("A117").PasteSpecial(-4122)
$ws.Range("H107").Copy()
$ws.Range("H117").PasteSpecial(-4122)
$ws.Range("A117").Value = "TOTAL"
$ws.Range("H117").Value = 1283.87

# Now fill in the actual data for the new block
$ws.Range("A112").Value = "Point 19"
$ws.Range("B112").Value = "ARM-3SF-GN-C"
$ws.Range("C112").Value = "Rem"
$ws.Range("D112").Value = "ARM,3ft Sgl. Fiberglass,Gain,Corr"
$ws.Range("E112").Value = "EA"
$ws.Range("F112").Value = 1
$ws.Range("H112").Value = 61.83

$ws.Range("A113").Value = "Point 19"
$ws.Range("B113").Value = "DEC-20AL-C"
$ws.Range("C113").Value = "Rem"
$ws.Range("D113").Value = "DEC,#4 - #2/0 AA,AL,AS,Corrosive"
$ws.Range("E113").Value = "EA"
$ws.Range("F113").Value = 6
$ws.Range("H113").Value = 570.9

$ws.Range("A114").Value = "Point 19"
$ws.Range("B114").Value = "POL-45-4"
$ws.Range("C114").Value = "Rem"
$ws.Range("D114").Value = "Pole,45ft,Class 4"
$ws.Range("E114").Value = "EA"
$ws.Range("F114").Value = 1
$ws.Range("H114").Value = 198.88

$ws.Range("A115").Value = "Point 19"
$ws.Range("B115").Value = "SAA-DI-2-C"
$ws.Range("C115").Value = "Rem"
$ws.Range("D115").Value = "SAA,Dead End Clamp I Bolt,2,Corr"
$ws.Range("E115").Value = "EA"
$ws.Range("F115").Value = 1
$ws.Range("H115").Value = 17.2

$ws.Range("A116").Value = "Point 15"
$ws.Range("B116").Value = "ANC-DHM-10-84-T1-C"
$ws.Range("C116").Value = "Inst"
$ws.Range("D116").Value = "ANC,Dbl Hlx Mach,10in,84in,TpEye 1in,Cor"
$ws.Range("E116").Value = "EA"
$ws.Range("F116").Value = 2
$ws.Range("H116").Value = 435.06
